$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-gender")

# Update is_active (column D) from TRUE to FALSE for the three "Others" rows
# Row 4: code=OTH, name=Others, lang_code=eng
# Row 7: code=OTH, name=الآخرين, lang_code=ara
# Row 10: code=OTH, name=Dautres, lang_code=fra
$ws.Range("D4").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("D10").Value = $false

# Update the active cell selection shown in the sheet view
$ws.Range("D12").Select()
